$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.916.42"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.874.82"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7403"
$ws.Range("E5").Value = "  -4.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.58"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3153"
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07161"
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.69"
$ws.Range("E10").Value = "  -3.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08401"
$ws.Range("E11").Value = "  -3.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7510"
$ws.Range("E12").Value = "  -2.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.421"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "1.853.65"
$ws.Range("E14").Value = "  -10.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.57"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "29.894.81"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.106"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("E18").Value = "  -2.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.17"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007813"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9985"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "2.121.68"
$ws.Range("E22").Value = "  -10.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.976"
$ws.Range("E23").Value = "  -2.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1549"
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.290"
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.16"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.62"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.036"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.489"
$ws.Range("E30").Value = "  +4.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.598"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.534"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.260"
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05326"
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.237"
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7549"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9960"
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.691"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01951"
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4508"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "1.111.84"
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.049"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.18"
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8558"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.03"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.657"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.102"
$ws.Range("E49").Value = "  +3.53%  "
$ws.Range("E50").Value = "  -2.57%  "
$ws.Range("D51").Value = "2.017.72"
$ws.Range("E51").Value = "  -9.90%  "
